# Update countries & provincias Spain
# Applies the updated COVID-19 country stats and the re-sorted country
# labels (rows whose "Casos totales" ranking changed) from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos
$ws.Range("B4").Value = 112533
$ws.Range("C4").Value = 8407
$ws.Range("D4").Value = 3219
$ws.Range("E4").Value = 107438
$ws.Range("F4").Value = 2666
$ws.Range("G4").Value = 180
$ws.Range("H4").Value = 1876

# Row 12: Suiza
$ws.Range("F12").Value = 280

# Row 19: Canada
$ws.Range("B19").Value = 4782
$ws.Range("C19").Value = 25
$ws.Range("E19").Value = 4373

# Row 20: Noruega
$ws.Range("B20").Value = 3973
$ws.Range("C20").Value = 202
$ws.Range("E20").Value = 3946

# Row 34: Rumania
$ws.Range("E34").Value = 1283
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 30

# Row 62: Barein
$ws.Range("B62").Value = 476
$ws.Range("C62").Value = 10
$ws.Range("D62").Value = 265
$ws.Range("E62").Value = 207

# Row 71: Bulgaria
$ws.Range("D71").Value = 11
$ws.Range("E71").Value = 296

# Row 83: Moldavia
$ws.Range("A83").Value = "Moldavia"
$ws.Range("B83").Value = 231
$ws.Range("C83").Value = 32
$ws.Range("E83").Value = 227
$ws.Range("F83").Value = 33
$ws.Range("H83").Value = 2

# Row 84: Kazajistan
$ws.Range("A84").Value = "Kazajistan"
$ws.Range("B84").Value = 228
$ws.Range("C84").Value = 78
$ws.Range("D84").Value = 16
$ws.Range("E84").Value = 211
$ws.Range("F84").Value = 0
$ws.Range("H84").Value = 1

# Row 85: Tunez
$ws.Range("A85").Value = "Tunez"
$ws.Range("B85").Value = 227
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 2
$ws.Range("E85").Value = 218
$ws.Range("F85").Value = 10
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 7

# Row 86: San Marino
$ws.Range("A86").Value = "San Marino"
$ws.Range("B86").Value = 223
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 6
$ws.Range("E86").Value = 196
$ws.Range("F86").Value = 15
$ws.Range("H86").Value = 21

# Row 87: Burkina Faso
$ws.Range("A87").Value = "Burkina Faso"
$ws.Range("B87").Value = 207
$ws.Range("C87").Value = 27
$ws.Range("D87").Value = 21
$ws.Range("E87").Value = 175
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 11

# Row 88: Albania
$ws.Range("A88").Value = "Albania"

# Row 89: Azerbaiyan
$ws.Range("A89").Value = "Azerbaiyan"
$ws.Range("B89").Value = 182
$ws.Range("C89").Value = 17
$ws.Range("D89").Value = 15
$ws.Range("E89").Value = 163
$ws.Range("F89").Value = 23
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 4

# Row 90: Republica de Chipre
$ws.Range("A90").Value = "Republica de Chipre"
$ws.Range("B90").Value = 179
$ws.Range("C90").Value = 17
$ws.Range("E90").Value = 159
$ws.Range("F90").Value = 3
$ws.Range("H90").Value = 5

# Row 91: Vietnam
$ws.Range("A91").Value = "Vietnam"
$ws.Range("B91").Value = 174
$ws.Range("C91").Value = 11
$ws.Range("D91").Value = 21
$ws.Range("E91").Value = 153
$ws.Range("H91").Value = 0

# Row 96: Ghana
$ws.Range("B96").Value = 141
$ws.Range("C96").Value = 4
$ws.Range("E96").Value = 134
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 5

# Row 113: Georgia
$ws.Range("A113").Value = "Georgia"
$ws.Range("B113").Value = 90
$ws.Range("C113").Value = 7
$ws.Range("D113").Value = 14
$ws.Range("E113").Value = 76
$ws.Range("F113").Value = 1
$ws.Range("H113").Value = 0

# Row 114: Nigeria
$ws.Range("A114").Value = "Nigeria"
$ws.Range("B114").Value = 89
$ws.Range("C114").Value = 19
$ws.Range("D114").Value = 3
$ws.Range("E114").Value = 85
$ws.Range("F114").Value = 0
$ws.Range("H114").Value = 1
